$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" column header - copy the bold/bordered header style from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Label values for rows 2-11 (first block, Iterations=100) and 12-21 (second block, Iterations=200)
# Controls (rows 2-6, 12-16) => 0 ; MDD (rows 7-11, 17-21) => 1
$labels = @(0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,1,1,1,1,1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}

# Updated Prediction/Error values (refit results)
$ws.Range("D2").Value = 0.6680081869367536
$ws.Range("E2").Value = 0.6680081869367536

$ws.Range("D8").Value = 0.6643967590047086
$ws.Range("E8").Value = 0.3356032409952914

$ws.Range("D11").Value = 0.7299772696453194
$ws.Range("E11").Value = 0.2700227303546806
